# Revert commit 6dc3ea4 ("Working on the ability to dump individual files, as
# well as fixing a major bug, on file 72"): the referenced shared strings for
# these 7 rows (columns D/E) had drifted onto duplicate string-table entries
# (identical text, just re-inserted with slightly different leading spacing).
# Restore each cell to its original text so the workbook again uses the single,
# de-duplicated copy of each string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D15").Value = 'Background
id="Par1">Previous studies have shown that Coronavirus Disease 2019 (COVID-19) patients with underlying comorbidities can have worse outcomes.

 However, the effect of hypertension on outcomes of COVID-19 patients remains unclear.


Research question
id="Par2">The aim of this study was to explore the effect of hypertension on the outcomes of patients with COVID-19 by using propensity score–matching (PSM) analysis.


Study design and methods
id="Par3">Participants enrolled in this study were patients with COVID-19 who had been hospitalized at the Central Hospital of Wuhan, China.

 Chronic comorbidities and laboratory and radiological data were reviewed; patient outcomes and lengths of stay were obtained from discharge records.

 We used the Cox proportional-hazard model (CPHM) to analyze the effect of hypertension on these patients’ outcomes and PSM analysis to further validate the abovementioned effect.


Results
id="Par4">A total of 226 patients with COVID-19 were enrolled in this study, of whom 176 survived and 50 died.

 The proportion of patients with hypertension among non-survivors was higher than that among survivors (26.70% vs.

 74.00%; P &lt; 0.001).

 Results obtained via CPHM showed that hypertension could increase risk of mortality in COVID-19 patients (hazard ratio 3.317; 95% CI [1.709–6.440]; P &lt; 0.001).

 Increased D-dimer levels and higher ratio of neutrophils to lymphocytes (N/L) were also found to increase these patients’ mortality risk.

 After matching on propensity score, we still came to similar conclusions.

 After we applied the same method in critically ill patients, we found that hypertension also increased risk of death in patients with severe COVID-19.
Conclusion
id="Par5">Hypertension, increased D-dimer and the ratio of neutrophil to lymphocyte increased mortality in patients with COVID-19, with hypertension in particular.


'

$ws.Range("E15").Value = '[Qing%Yang%NULL%0,  Ying%Zhou%NULL%0,  Xinrong%Wang%NULL%0,  Shan%Gao%NULL%0,  Yang%Xiao%NULL%0,  Weiming%Zhang%NULL%0,  Yi%Hu%huyizxyy@163.com%0,  Yafei%Wang%wyf_527@163.com%0]'

$ws.Range("D17").Value = 'Background
id="Par1">Since December 2019, the outbreak of COVID-19 caused a large number of hospital admissions in China.

 Many patients with COVID-19 have symptoms of acute respiratory distress syndrome, even are in danger of death.

 This is the first study to evaluate dynamic changes of D-Dimer and Neutrophil-Lymphocyte Count Ratio (NLR) as a prognostic utility in patients with COVID-19 for clinical use.


Methods
id="Par2">In a retrospective study, we collected data from 349 hospitalized patients who diagnosed as the infection of the COVID-19 in Wuhan Pulmonary Hospital.

 We used ROC curves and Cox regression analysis to explore critical value (optimal cut-off point associated with Youden index) and prognostic role of dynamic changes of D-Dimer and NLR.


Results
id="Par3">Three hundred forty-nine participants were enrolled in this study and the mortality rate of the patients with laboratory diagnosed COVID-19 was 14.9%.

 The initial and peak value of D-Dimer and NLR in deceased patients were higher statistically compared with survivors (P &lt; 0.001).

 There was a more significant upward trend of D-Dimer and NLR during hospitalization in the deceased patients, initial D-Dimer and NLR were lower than the peak tests (MD) -25.23, 95% CI: − 31.81- -18.64, P &lt; 0.001; (MD) -43.73, 95% CI:-59.28- -31.17, P &lt; 0.001. The test showed a stronger correlation between hospitalization days, PCT and peak D-Dimer than initial D-Dimer.

 The areas under the ROC curves of peak D-Dimer and peak NLR tests were higher than the initial tests (0.94(95%CI: 0.90–0.98) vs.

 0.80 (95% CI: 0.73–0.87); 0.93 (95%CI:0.90–0.96) vs.

 0.86 (95%CI:0.82–0.91).

 The critical value of initial D-Dimer, peak D-Dimer, initial NLR and peak NLR was 0.73 mg/L, 3.78 mg/L,7.13 and 14.31 respectively.

 35 (10.03%) patients were intubated.

 In the intubated patients, initial and peak D-Dimer and NLR were much higher than non-intubated patients (P &lt; 0.001).

 The critical value of initial D-Dimer, peak D-Dimer, initial NLR and peak NLR in prognosticate of intubation was 0.73 mg/L, 12.75 mg/L,7.28 and 27.55. The multivariable Cox regression analysis showed that age (HR 1.04, 95% CI 1.00–1.07, P = 0.01), the peak D-Dimer (HR 1.03, 95% CI 1.01–1.04, P &lt; 0.001) were prognostic factors for COVID-19 patients’ death.


Conclusions
id="Par4">To dynamically observe the ratio of D-Dimer and NLR was more valuable during the prognosis of COVID-19. The rising trend in D-Dimer and NLR, or the test results higher than the critical values may indicate a risk of death for participants with COVID-19.
'

$ws.Range("E17").Value = '[Wenjing%Ye%yewenjing@xinhuamed.com.cn%0,  Guoxi%Chen%4799082@qq.com%0,  Xiaopan%Li%xiaopanli0224@126.com%0,  Xing%Lan%474137452@qq.com%0,  Chen%Ji%c.ji.3@warwick.ac.uk%0,  Min%Hou%min-hou710@hotmail.com%0,  Di%Zhang%zhangdizhangdi1234@163.com%0,  Guangwang%Zeng%785663609@qq.com%0,  Yaling%Wang%953822402@qq.com%0,  Cheng%Xu%xucheng112358@126.com%0,  Weiwei%Lu%luweiwei100@126.com%0,  Ruolin%Cui%crlqwerty@163.com%0,  Yuyang%Cai%caiyuyang@sjtu.edu.cn%0,  Hai%Huang%1220775601@qq.com%0,  Ling%Yang%yangling01@xinhuamed.com.cn%0]'

$ws.Range("D26").Value = 'Background
id="Par1">Data from the UK COVID-19 outbreak are emerging, and there are ongoing concerns about a disproportionate effect on ethnic minorities.

 There is very limited information on COVID-19 in the over-80s, and the rates of hospital-onset infections are unknown.


Methods
id="Par2">This was a retrospective cohort study from electronic case records of the first 450 patients admitted to our hospital with PCR-confirmed COVID-19, 77% of the total inpatient caseload to date.

 Demographic, clinical and biochemical data were extracted.

 The primary endpoint was death during the index hospital admission.

 The characteristics of all patients, those over 80 years of age and those with hospital-onset COVID-19 were examined.


Results
id="Par3">The median (IQR) age was 72 (56, 83), with 150 (33%) over 80 years old and 60% male.

 Presenting clinical and biochemical features were consistent with those reported elsewhere.

 The ethnic breakdown of patients admitted was similar to that of our underlying local population.

 Inpatient mortality was high at 38%.


Conclusions
id="Par5">Inpatient mortality was high, especially among the over-80s, who are more likely to present atypically.

 The ethnic composition of our caseload was similar to the underlying population.

 While a significant number of patients acquired COVID-19 while already in hospital, their outcomes were no worse.


'

$ws.Range("E26").Value = '[Simon E.%Brill%simon.brill@nhs.net%0,  Hannah C.%Jarvis%NULL%0,  Hannah C.%Jarvis%NULL%0,  Ezgi%Ozcan%NULL%0,  Thomas L. P.%Burns%NULL%0,  Rabia A.%Warraich%NULL%0,  Lisa J.%Amani%NULL%0,  Amina%Jaffer%NULL%0,  Stephanie%Paget%NULL%0,  Anand%Sivaramakrishnan%NULL%0,  Dean D.%Creer%NULL%0]'

$ws.Range("D43").Value = 'Introduction
id="Par1">In the current study, we evaluated factors that increase the coronavirus disease (COVID-19) patient death rate by analyzing the data from two cohort hospitals.

 In addition, we studied whether underlying neurological diseases are risk factors for death.


Methods
id="Par2">In this retrospective cohort study, we included 103 adult inpatients (aged ≥ 18 years).

 We evaluated differences in demographic data between surviving and non-surviving COVID-19 patients.


Results
id="Par3">In a multivariate logistic analysis, age and the presence of chronic lung disease and Alzheimer’s dementia (AD) were the only significant parameters for predicting COVID-19 non-survival (p &lt; 0.05).

 However, hypertension, coronary vascular disease, dyslipidemia, chronic kidney disease, diabetes, and history of taking angiotensin II receptor blockers (ARBs) or angiotensin-converting enzyme (ACE) inhibitors, as well as nonsteroidal anti-inflammatory drugs (NSAIDs), were not significantly associated with the death of COVID-19 patients.

 The optimal cutoff value obtained from the maximum Youden index was 70 (sensitivity, 80.77%; specificity, 61.04%), and the odds ratio of non-survival increased 1.055 fold for every year of age.


Conclusions
id="Par4">Clinicians should closely monitor and manage the symptoms of COVID-19 patients who are over the age of 70 years or have chronic lung disease or AD.


'

$ws.Range("E43").Value = '[Jong-moon%Hwang%NULL%0,  Ju-Hyun%Kim%NULL%0,  Jin-Sung%Park%NULL%0,  Min Cheol%Chang%wheel633@ynu.ac.kr%0,  Donghwi%Park%bdome@hanmail.net%0]'

$ws.Range("D54").Value = 'Background
id="Par1">There are currently rare satisfactory markers for predicting the death of patients with coronavirus disease 2019 (COVID-19).

 The aim of this study is to establish a model based on the combination of serum cytokines and lymphocyte subsets for predicting the prognosis of the disease.


Methods
id="Par2">A total of 739 participants with COVID-19 were enrolled at Tongji Hospital from February to April 2020 and classified into fatal (n = 51) and survived (n = 688) groups according to the patient’s outcome.

 Cytokine profile and lymphocyte subset analysis was performed simultaneously.


Results
id="Par3">The fatal patients exhibited a significant lower number of lymphocytes including B cells, CD4+ T cells, CD8+ T cells, and NK cells and remarkably higher concentrations of cytokines including interleukin-2 receptor, interleukin-6, interleukin-8, and tumor necrosis factor-α on admission compared with the survived subjects.

 A model based on the combination of interleukin-8 and the numbers of CD4+ T cells and NK cells showed a good performance in predicting the death of patients with COVID-19. When the threshold of 0.075 was used, the sensitivity and specificity of the prediction model were 90.20% and 90.26%, respectively.

 Meanwhile, interleukin-8 was found to have a potential value in predicting the length of hospital stay until death.


Conclusions
id="Par4">Significant increase of cytokines and decrease of lymphocyte subsets are found positively correlated with in-hospital death.

 A model based on the combination of three markers provides an attractive approach to predict the prognosis of COVID-19.
Electronic supplementary material
The online version of this article (10.1007/s10875-020-00821-7) contains supplementary material, which is available to authorized users.


'

$ws.Range("E54").Value = '[Ying%Luo%NULL%0,  Liyan%Mao%NULL%0,  Xu%Yuan%NULL%0,  Ying%Xue%NULL%0,  Qun%Lin%NULL%0,  Guoxing%Tang%NULL%0,  Huijuan%Song%NULL%0,  Feng%Wang%fengwang@tjh.tjmu.edu.cn%0,  Ziyong%Sun%zysun@tjh.tjmu.edu.cn%0]'

$ws.Range("D57").Value = 'Background
id="Par1">New York City emerged as an epicenter of the coronavirus disease 2019 (COVID-19) pandemic.


Objective
id="Par2">To describe the clinical characteristics and risk factors associated with mortality in a large patient population in the USA.


Design
id="Par3">Retrospective cohort study.


Participants
id="Par4">6493 patients who had laboratory-confirmed COVID-19 with clinical outcomes between March 13 and April 17, 2020, who were seen in one of the 8 hospitals and/or over 400 ambulatory practices in the New York City metropolitan area
Main Measures
id="Par5">Clinical characteristics and risk factors associated with in-hospital mortality.


Key Results
id="Par6">A total of 858 of 6493 (13.2%) patients in our total cohort died: 52/2785 (1.9%) ambulatory patients and 806/3708 (21.7%) hospitalized patients.

 Cox proportional hazard regression modeling showed an increased risk of in-hospital mortality associated with age older than 50 years (hazard ratio [HR] 2.34, CI 1.47–3.71), systolic blood pressure less than 90 mmHg (HR 1.38, CI 1.06–1.80), a respiratory rate greater than 24 per min (HR 1.43, CI 1.13–1.83), peripheral oxygen saturation less than 92% (HR 2.12, CI 1.56–2.88), estimated glomerular filtration rate less than 60 mL/min/1.73m2 (HR 1.80, CI 1.60–2.02), IL-6 greater than 100 pg/mL (HR 1.50, CI 1.12–2.03), D-dimer greater than 2 mcg/mL (HR 1.19, CI 1.02–1.39), and troponin greater than 0.03 ng/mL (HR 1.40, CI 1.23–1.62).

 Decreased risk of in-hospital mortality was associated with female sex (HR 0.84, CI 0.77–0.90), African American race (HR 0.78 CI 0.65–0.95), and hydroxychloroquine use (HR 0.53, CI 0.41–0.67).


Conclusions
id="Par7">Among patients with COVID-19, older age, male sex, hypotension, tachypnea, hypoxia, impaired renal function, elevated D-dimer, and elevated troponin were associated with increased in-hospital mortality and hydroxychloroquine use was associated with decreased in-hospital mortality.


Electronic supplementary material
The online version of this article (10.1007/s11606-020-05983-z) contains supplementary material, which is available to authorized users.


'

$ws.Range("E57").Value = '[Takahisa%Mikami%NULL%0,  Hirotaka%Miyashita%NULL%0,  Takayuki%Yamada%NULL%0,  Matthew%Harrington%NULL%0,  Daniel%Steinberg%NULL%0,  Andrew%Dunn%NULL%0,  Evan%Siau%Evan.Siau@mountsinai.org%0]'

$ws.Range("D58").Value = 'Background
id="Par1">The purpose of this study is to report the clinical features and outcomes of Black/African American (AA) and Latino Hispanic patients with Coronavirus disease 2019 (COVID-19) hospitalized in an inter-city hospital in the state of New Jersey.


Methods
id="Par2">This is a retrospective cohort study of AA and Latino Hispanic patients with COVID-19 admitted to a 665-bed quaternary care, teaching hospital located in Newark, New Jersey.

 The study included patients who had completed hospitalization between March 10, 2020, and April 10, 2020. We reviewed demographics, socioeconomic variables and incidence of in-hospital mortality and morbidity.

 Logistic regression was used to identify predictor of in-hospital death.


Results
id="Par3">Out of 416 patients, 251 (60%) had completed hospitalization as of April 10, 2020. The incidence of In-hospital mortality was 38.6% (n = 97).

 Most common symptoms at initial presentation were dyspnea 39% (n = 162) followed by cough 38%(n = 156) and fever 34% (n = 143).

 Patients were in the highest quartile for population’s density, number of housing units and disproportionately fell into the lowest median income quartile for the state of New Jersey.

 The incidence of septic shock, acute kidney injury (AKI) requiring hemodialysis and admission to an intensive care unit (ICU) was 24% (n = 59), 21% (n = 52), 33% (n = 82) respectively.

 Independent predictors of in-hospital mortality were older age, lower serum Hemoglobin &lt; 10 mg/dl, elevated serum Ferritin and Creatinine phosphokinase levels &gt; 1200 U/L and &gt; 1000 U/L.


Conclusions
id="Par4">Findings from an inter-city hospital’s experience with COVID-19 among underserved minority populations showed that, more than one of every three patients were at risk for in-hospital death or morbidity.

 Older age and elevated inflammatory markers at presentation were associated with in-hospital death.


'

$ws.Range("E58").Value = '[Alexis K.%Okoh%alexis.okoh@rwjbh.org%0,  Christoph%Sossou%NULL%0,  Christoph%Sossou%NULL%0,  Neha S.%Dangayach%NULL%0,  Sherin%Meledathu%NULL%0,  Oluwakemi%Phillips%NULL%0,  Corinne%Raczek%NULL%0,  Michael%Patti%NULL%0,  Nathan%Kang%NULL%0,  Sameer A.%Hirji%NULL%0,  Charles%Cathcart%NULL%0,  Christian%Engell%NULL%0,  Marc%Cohen%NULL%0,  Sandhya%Nagarakanti%NULL%0,  Eliahu%Bishburg%NULL%0,  Harpreet S.%Grewal%NULL%0]'
